# Corrects methodological issue with ban on CCUS retrofits; updates input data

$wb = $excel.ActiveWorkbook

# --- BBNPPTY sheet: extend the "banned" flag (1) for CCUS-retrofit fuel types
#     through 2024-2027 (previously only banned from 2028 onward) ---
$wsData = $wb.Worksheets.Item("BBNPPTY")

$wsData.Range("E19:H19").Value = 1   # hard coal w CCS
$wsData.Range("E20:H20").Value = 1   # natural gas combined cycle w CCS
$wsData.Range("E21:H21").Value = 1   # biomass w CCS
$wsData.Range("E22:H22").Value = 1   # lignite w CCS

# leave the selection positioned over the data that was just updated
$wsData.Range("E19:H22").Select() | Out-Null

# --- About sheet: document the methodology correction with two new note rows ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A15").Value = "We also assume no new coal with CCS can be built prior to 2028 given the state of the technology"
$wsAbout.Range("A16").Value = "and the construction time for new or modified plants."

# leave the workbook showing the About sheet, selection resting below the new notes
$wsAbout.Activate() | Out-Null
$wsAbout.Range("B22").Select() | Out-Null
